{"js": "// Remove the sentence \"It also decays epsilon in order to get better\n// convergence speed.\" (including the line break that separated it from the\n// following \"At every checkpoint ...\" sentence), per the commit message:\n// \"removed the statement that says that epsilon is reduced\".\n//\n// Before: \"...go to the next one. It also decays epsilon in order to get\n//          better convergence speed.<line break>At every checkpoint...\"\n// After : \"...go to the next one. At every checkpoint...\"\n\nconst searchText =\n  \" It also decays epsilon in order to get better convergence speed.\" +\n  \"\\u000bAt every checkpoint\";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found in document body.\");\n}\n\n// Replace the matched range (sentence + following line break + the start of\n// the next sentence) with just \" At every checkpoint\" so the two sentences\n// end up on the same line, separated by a single space.\nresults.items[0].insertText(\" At every checkpoint\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Remove the sentence \"It also decays epsilon in order to get better\n# convergence speed.\" (including the line break that separated it from the\n# following \"At every checkpoint ...\" sentence), per the commit message:\n# \"removed the statement that says that epsilon is reduced\".\n#\n# Before: \"...go to the next one. It also decays epsilon in order to get\n#          better convergence speed.<line break>At every checkpoint...\"\n# After : \"...go to the next one. At every checkpoint...\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# \"^l\" is Word's Find/Replace wildcard for a manual line break (vertical tab),\n# so this pattern spans the sentence to remove plus the line break right\n# after it, up to the start of the next sentence.\n$find.Text = \" It also decays epsilon in order to get better convergence speed.^lAt every checkpoint\"\n$find.Replacement.Text = \" At every checkpoint\"\n\n$find.Forward = $true\n$find.Wrap = [Microsoft.Office.Interop.Word.WdFindWrap]::wdFindStop\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n$find.Execute([Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Type]::Missing, [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceOne) | Out-Null\n"}
